$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q4" sheet right after "总计", cloning the formatting
#    of the existing "2022-Q3" sheet (same header layout / styles) and then
#    replacing its data with the 2022-Q4 figures.
# ---------------------------------------------------------------------------
$refSheet   = $wb.Worksheets.Item("2022-Q3")
$totalSheet = $wb.Worksheets.Item(1)

$refSheet.Copy($null, $totalSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The template (2022-Q3) carried 3 data rows - drop the extra ones, keep row 2.
$newSheet.Range("A3:H4").Delete()

# Fill in the 2022-Q4 fund data (row 2).
$newSheet.Range("B2").Value = "'320022"
$newSheet.Range("C2").Value = "诺安研究精选股票"
$newSheet.Range("D2").Value = "'6.42"
$newSheet.Range("E2").Value = "'92.87"
$newSheet.Range("F2").Value = "'2.20"
$newSheet.Range("G2").Value = "'0.1412"
$newSheet.Range("H2").Value = 6
# Drop the "quote prefix" formatting that typing numeric-looking text leaves
# behind, so the cells end up unstyled like the rest of the sheet.
$newSheet.Range("B2:G2").ClearFormats()

# Restore the originally active tab (2020-Q4) as the selected sheet.
$wb.Worksheets.Item("2020-Q4").Activate()

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q4 above
#    the existing 2022-Q3 row, shifting everything else down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.14

# A column keeps the bordered/bold style used by every other row - copy it
# from the row below (which holds the old 2022-Q3 row) onto the new row.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
